# Update report: add 4 new "Đơn sale chính" rows, recompute the totals row,
# and add the "Đơn thu nợ" and "Lương" sheets (report_ca_nhan update).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Columns C (ngày thực hiện) hold dd-mm-yyyy text, not real dates. Pre-format
# as Text so Excel doesn't auto-convert the literal strings to date serials.
$ws1.Range("C3:C6").NumberFormat = "@"

# --- Sheet 1: "Đơn sale chính" ---------------------------------------------
# Row 3 used to be the "Tổng" (totals) row; it becomes a normal data row and
# three more data rows are appended, with a fresh totals row at row 7.

$ws1.Range("A3").Value = "HD-LUXURY"
$ws1.Range("B3").Value = 519
$ws1.Range("C3").Value = "07-05-2024"
$ws1.Range("D3").Value = "SÓC TRĂNG"
$ws1.Range("E3").Value = "na quy"
$ws1.Range("F3").Value = "CTV"
$ws1.Range("G3").Value = "Tiêm Filler"
$ws1.Range("H3").Value = "Thạch Hoàng Nhân"
$ws1.Range("I3").Value = 21500000
$ws1.Range("J3").Value = 0
$ws1.Range("K3").Value = 0
$ws1.Range("L3").Value = 21500000
$ws1.Range("M3").Value = 21500000
$ws1.Range("N3").Value = 0
$ws1.Range("O3").Value = 21500000
$ws1.Range("P3").Value = 0
$ws1.Range("Q3").Value = "Nguyễn Hoàng Yến Quyên"
$ws1.Range("R3").Value = 0
$ws1.Range("S3").Value = "Kha Như Huỳnh "
$ws1.Range("T3").Value = 0
$ws1.Range("U3").Value = 50000
$ws1.Range("V3").Value = 0

$ws1.Range("A4").Value = "HD-LUXURY"
$ws1.Range("B4").Value = 520
$ws1.Range("C4").Value = "07-05-2024"
$ws1.Range("D4").Value = "SÓC TRĂNG"
$ws1.Range("E4").Value = "đa ni "
$ws1.Range("F4").Value = "CTV"
$ws1.Range("G4").Value = "Tiêm Filler"
$ws1.Range("H4").Value = "Thạch Hoàng Nhân"
$ws1.Range("I4").Value = 11000000
$ws1.Range("J4").Value = 0
$ws1.Range("K4").Value = 0
$ws1.Range("L4").Value = 11000000
$ws1.Range("M4").Value = 11000000
$ws1.Range("N4").Value = 0
$ws1.Range("O4").Value = 11000000
$ws1.Range("P4").Value = 0
$ws1.Range("Q4").Value = 0
$ws1.Range("R4").Value = 0
$ws1.Range("S4").Value = "Kha Như Huỳnh "
$ws1.Range("T4").Value = 0
$ws1.Range("U4").Value = 50000
$ws1.Range("V4").Value = 0

$ws1.Range("A5").Value = "HD-LUXURY"
$ws1.Range("B5").Value = 521
$ws1.Range("C5").Value = "07-05-2024"
$ws1.Range("D5").Value = "SÓC TRĂNG"
$ws1.Range("E5").Value = "đa ni "
$ws1.Range("F5").Value = "CTV"
$ws1.Range("G5").Value = "Phun môi"
$ws1.Range("H5").Value = "Thạch Hoàng Nhân"
$ws1.Range("I5").Value = 5500000
$ws1.Range("J5").Value = 0
$ws1.Range("K5").Value = 0
$ws1.Range("L5").Value = 5500000
$ws1.Range("M5").Value = 5500000
$ws1.Range("N5").Value = 0
$ws1.Range("O5").Value = 5500000
$ws1.Range("P5").Value = 0
$ws1.Range("Q5").Value = "Bác Sĩ Ngoài"
$ws1.Range("R5").Value = 0
$ws1.Range("S5").Value = 0
$ws1.Range("T5").Value = 0
$ws1.Range("U5").Value = 0
$ws1.Range("V5").Value = 0

$ws1.Range("A6").Value = "HD-LUXURY"
$ws1.Range("B6").Value = 522
$ws1.Range("C6").Value = "07-05-2024"
$ws1.Range("D6").Value = "SÓC TRĂNG"
$ws1.Range("E6").Value = "ngọc linh "
$ws1.Range("F6").Value = "CTV"
$ws1.Range("G6").Value = "Làm má Lúm"
$ws1.Range("H6").Value = "Thạch Hoàng Nhân"
$ws1.Range("I6").Value = 3000000
$ws1.Range("J6").Value = 0
$ws1.Range("K6").Value = 0
$ws1.Range("L6").Value = 3000000
$ws1.Range("M6").Value = 3000000
$ws1.Range("N6").Value = 0
$ws1.Range("O6").Value = 3000000
$ws1.Range("P6").Value = 0
$ws1.Range("Q6").Value = 0
$ws1.Range("R6").Value = 0
$ws1.Range("S6").Value = "Kha Như Huỳnh "
$ws1.Range("T6").Value = 0
$ws1.Range("U6").Value = 50000
$ws1.Range("V6").Value = 0

$ws1.Range("A7").Value = "Tổng"
$ws1.Range("B7").Value = 5
$ws1.Range("C7").Value = ""
$ws1.Range("D7").Value = ""
$ws1.Range("E7").Value = ""
$ws1.Range("F7").Value = ""
$ws1.Range("G7").Value = ""
$ws1.Range("H7").Value = ""
$ws1.Range("I7").Value = 66000000
$ws1.Range("J7").Value = ""
$ws1.Range("K7").Value = 0
$ws1.Range("L7").Value = 66000000
$ws1.Range("M7").Value = 66000000
$ws1.Range("N7").Value = 0
$ws1.Range("O7").Value = 66000000
$ws1.Range("P7").Value = 0
$ws1.Range("Q7").Value = ""
$ws1.Range("R7").Value = ""
$ws1.Range("S7").Value = ""
$ws1.Range("T7").Value = ""
$ws1.Range("U7").Value = 200000
$ws1.Range("V7").Value = 0

# --- Sheet 2: "Đơn thu nợ" ---------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Đơn thu nợ"

$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

$ws2.Range("G2:G3").NumberFormat = "@"

$ws2.Range("A1").Value = "Tiền tố"
$ws2.Range("B1").Value = "Mã đơn thu nợ"
$ws2.Range("C1").Value = "Đơn nợ"
$ws2.Range("D1").Value = "Cơ sở"
$ws2.Range("E1").Value = "Lượng thu"
$ws2.Range("F1").Value = "Sale"
$ws2.Range("G1").Value = "Ngày thu"

$ws2.Range("A2").Value = "TN"
$ws2.Range("B2").Value = 137
$ws2.Range("C2").Value = "HD-LUXURY-262"
$ws2.Range("D2").Value = "SÓC TRĂNG"
$ws2.Range("E2").Value = 10000000
$ws2.Range("F2").Value = "Thạch Hoàng Nhân"
$ws2.Range("G2").Value = "07-05-2024"

$ws2.Range("A3").Value = "TN"
$ws2.Range("B3").Value = 139
$ws2.Range("C3").Value = "HD-LUXURY-356"
$ws2.Range("D3").Value = "SÓC TRĂNG"
$ws2.Range("E3").Value = 2000000
$ws2.Range("F3").Value = "Thạch Hoàng Nhân"
$ws2.Range("G3").Value = "07-05-2024"

$ws2.Range("A4").Value = "Tổng"
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = ""
$ws2.Range("D4").Value = ""
$ws2.Range("E4").Value = 12000000
$ws2.Range("F4").Value = ""
$ws2.Range("G4").Value = ""

# --- Sheet 3: "Lương" --------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Lương"

$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

$ws1.Activate()
